$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-02-28"

# Update the row label text for February
$ws.Range("A3").Value = "February (through 02-28)"

# Update February row (row 3) values
$ws.Range("B3").Value = 11
$ws.Range("D3").Value = 56
$ws.Range("F3").Value = 30
$ws.Range("H3").Value = 125
$ws.Range("I3").Value = 141

# Update Total row (row 4) values
$ws.Range("B4").Value = 37
$ws.Range("D4").Value = 131
$ws.Range("F4").Value = 79
$ws.Range("H4").Value = 342
$ws.Range("I4").Value = 300
